$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at row 382 (pushes the old rows 382..484 down
# to 383..485, matching the diff which shows every record shifting down by
# one row and a fresh record appearing at the top of the block).
$ws.Rows("382:382").Insert()

$newRow = 382
$ws.Cells.Item($newRow, 1).Value = 8
$ws.Cells.Item($newRow, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item($newRow, 3).Value = "Coquimbo"
$ws.Cells.Item($newRow, 4).Value = 44642
$ws.Cells.Item($newRow, 5).Value = 4
$ws.Cells.Item($newRow, 6).Value = 100112043
$ws.Cells.Item($newRow, 7).Value = "Pepino ensalada"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 500
$ws.Cells.Item($newRow, 11).Value = 18000
$ws.Cells.Item($newRow, 12).Value = 19000
$ws.Cells.Item($newRow, 13).Value = 18500
$ws.Cells.Item($newRow, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item($newRow, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($newRow, 16).Value = 308
$ws.Cells.Item($newRow, 17).Value = 60
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
